$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching style of existing headers (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add new data value in H2
$ws.Range("H2").Value = 0
